$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (JHeater, JRTD connector) was manually swapped from the vertical
# PicoBlade 53398 header to the right-angle PicoBlade 53261 header.
$ws.Range("F13").Value = 532610871
$ws.Range("G13").Value = "Connector Header Surface Mount, Right Angle 8 position 0.049"" (1.25mm)"
$ws.Range("H13").Value = "PicoBlade 53261"
$ws.Range("I13").Value = 532610871

# Reflect the row growing taller to fit the new, longer description text.
$ws.Rows.Item(13).RowHeight = 63.75

# Move the active selection like the author left it after editing.
$ws.Range("H14").Select()
